$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update "Days remaining" values for two trials (B6: REJOICE, B8: REMASTER (CLOU))
$ws.Range("B6").Value = 1
$ws.Range("B8").Value = 21
